$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    $idx = 0
    $found = -1
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            $found = $idx
        }
    }
    return $found
}

# ---------------------------------------------------------------
# 1. Merge the two runs "...possibl" + "e." into a single run that
#    reads "...possible."
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "hot summer day for as long as possible.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "hot summer day for as long as possible.", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Fix "to long" -> "too long"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "to long", $true, $false, $false, $false, $false, $true, 1, $false,
    "too long", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Insert a brand-new paragraph "ANYONE with free time" right
#    after the "Everyone who..." paragraph.
# ---------------------------------------------------------------
$audienceIdx = Get-ParaIndexByText $d "Everyone who has delt with the frustration of keeping something out of the fridge too long"
$audiencePara = $d.Paragraphs.Item($audienceIdx)
$audiencePara.Range.InsertAfter([char]13 + "ANYONE with free time")

# ---------------------------------------------------------------
# 4. Fill the (previously empty) paragraph right after the
#    "Key features & Mechanics" heading with the mechanics blurb.
# ---------------------------------------------------------------
$keyFeaturesIdx = Get-ParaIndexByText $d "Key features & Mechanics"
$mechanicsPara = $d.Paragraphs.Item($keyFeaturesIdx + 1)
$mechanicsPara.Range.Text = "The player must use items to either stop a stick of butter from melting or cause it to melt. The player will gain coins the longer the butter is kept alive. These coins will allow the player to purchase items to make keeping the butter alive easier. These items can include a fan, wind, ice storm, snow, lower temperature, heat wave, hot breeze, heater, oven, or fire.  "

Write-Output "Done with text edits"
